$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "1.004",
# "24.743.51") that must stay plain text, exactly like the source
# workbook's inline-string cells. Writing such a string straight into
# .Value lets Excel "helpfully" reinterpret it as a number (dropping
# trailing zeros, switching to scientific notation, etc.), so the
# column is temporarily forced to Text format for the duration of the
# writes and then restored to the workbook's default (unstyled) look.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Cell updates derived from the authoritative diff (cell ref -> new value).
$ws.Range('D2').Value = '24.743.51'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '1.697.98'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '317.47'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '0.3929'
$ws.Range('E7').Value = '  -0.58%  '
$ws.Range('D8').Value = '0.4041'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = '1.506'
$ws.Range('E9').Value = '  -2.94%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '54.13'
$ws.Range('E10').Value = '  -5.43%  '
$ws.Range('B11').Value = 'BinanceUSD'
$ws.Range('C11').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D11').Value = '1.004'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '0.08891'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').Value = '7.261'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = '23.37'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').Value = '8.021'
$ws.Range('E15').Value = '  +4.68%  '
$ws.Range('D16').Value = '0.00001320'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '1.687.29'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').Value = '100.14'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('D19').Value = '0.07032'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').Value = '19.62'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').Value = '7.008'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '14.48'
$ws.Range('E23').Value = '  +2.52%  '
$ws.Range('D24').Value = '24.744.24'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').Value = '3.246'
$ws.Range('E25').Value = '  +9.34%  '
$ws.Range('D26').Value = '2.355'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('D27').Value = '22.82'
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('D28').Value = '161.45'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').Value = '136.86'
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('D30').Value = '5.182'
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('D31').Value = '7.760'
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '1.081'
$ws.Range('E32').Value = '  -3.20%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.08736'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('D34').Value = '7.205'
$ws.Range('E34').Value = '  -4.89%  '
$ws.Range('D35').Value = '11.35'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '1.961'
$ws.Range('E36').Value = '  -1.58%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').Value = '0.2749'
$ws.Range('E37').Value = '  -0.67%  '
$ws.Range('D38').Value = '14.38'
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('D39').Value = '0.09199'
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('D40').Value = '0.02736'
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('D41').Value = '1.466'
$ws.Range('E41').Value = '  -0.51%  '
$ws.Range('D42').Value = '0.7688'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').Value = '15.97'
$ws.Range('E43').Value = '  +1.35%  '
$ws.Range('D44').Value = '0.7174'
$ws.Range('E44').Value = '  -1.54%  '
$ws.Range('D45').Value = '2.561'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').Value = '4.221'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').Value = '1.002'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = '140.43'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').Value = '1.320'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('D50').Value = '90.66'
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('D51').Value = '0.07986'
$ws.Range('E51').Value = '  -0.39%  '

# Restore the Price column to the workbook's default (no explicit) style
# now that the text values are safely stored, so the cells end up
# formatted exactly as they were before (no stray number format).
$priceRange.Style = "Normal"
